$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All target cells are plain text
# (inline-string-typed Price / Volume(1h) figures), so force the Text number
# format before assigning the value. That stops Excel from auto-coercing
# "323.27" into a float or "-2.00%" into a computed percentage, which would
# silently change the stored precision/representation.
$updates = @{
    "D2" = "323.27"
    "E2" = "-2.00%"
    "D3" = "39.39"
    "E3" = "-1.60%"
    "D4" = "5.703"
    "E4" = "8.47%"
    "D5" = "0.07984"
    "E5" = "-1.41%"
    "D6" = "8.601"
    "E6" = "-0.46%"
    "D7" = "1.977"
    "E7" = "2.76%"
    "D8" = "2.945"
    "E8" = "-0.42%"
    "D9" = "0.9265"
    "E9" = "-0.80%"
    "D10" = "0.1257"
    "E10" = "-5.42%"
    "D11" = "0.1965"
    "E11" = "-0.11%"
    "D12" = "8.714"
    "E12" = "24.85%"
    "D13" = "0.09181"
    "E13" = "0.68%"
    "D14" = "0.03663"
    "E14" = "3.62%"
    "E15" = "9.54%"
    "E16" = "-2.07%"
    "D17" = "0.006237"
    "E17" = "5.19%"
    "D18" = "3.350"
    "D19" = "4.536"
    "E19" = "0.06%"
    "D20" = "0.3534"
    "E20" = "0.57%"
    "E21" = "3.21%"
    "D22" = "0.2453"
    "E22" = "-4.17%"
    "D23" = "0.04417"
    "E23" = "-0.03%"
    "D24" = "0.001266"
    "E24" = "3.67%"
    "D25" = "0.004552"
    "E25" = "5.37%"
    "D26" = "0.0001151"
    "E26" = "-3.15%"
    "D39" = "0.02503"
    "E39" = "-0.15%"
    "D40" = "0.05376"
    "E40" = "3.99%"
    "D41" = "0.007455"
    "E41" = "-3.17%"
    "D42" = "0.009562"
    "E42" = "3.73%"
    "D43" = "0.1403"
    "E43" = "-1.68%"
    "D44" = "0.002119"
    "E44" = "-1.80%"
    "D45" = "0.01088"
    "E45" = "-0.96%"
    "D46" = "0.00006780"
    "E46" = "1.85%"
    "E47" = "0.34%"
    "D48" = "0.002974"
    "E48" = "-11.09%"
    "D50" = "0.00002106"
    "E50" = "0.34%"
    "D51" = "0.0002005"
    "E51" = "0.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
